$d = $word.ActiveDocument

# --- Change 1: merge the 4 runs describing the islands/void text into a single run. ---
# The visible text is unchanged; Word's Find & Replace naturally collapses the
# matched range into one run (taking the formatting of the first run in the match).
$mergedText = " generatie die daar is/ heeft gewoond. Je bent op een stel eilanden ( het zijn geen vliegende eilanden of iets zoals dat, het zijn platformen en de grond eronder gaat heel diep de void in ). En nadat je de tutorial heb complete val je op een groot eiland en dat is de begin area, vanaf daar ga je parkouren en enemies bevechten op kleinere eilanden en dit doe je tot je weer bij een groter eiland komt waar een ruïne op staat en in die ruïne moet je een bepaalde puzzel oplossen. Je hebt 2 ruïnes en dus ook 2 puzzels, beide puzzels zijn hetzelfde idee maar wel een beetje anders, bij beide puzzels moet je een bal(len) van de startplek naar een knop brengen ******SIMON TYPE HIER JOUW GEWELDIGE IDEE******. Als je dit doet bij de eerste ruïne en je hebt de puzzel complete dan spawned er een item/ power up voor die je nodig hebt in de 2"
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# --- Change 2: after the "Game Mechanics:" paragraph, add a new paragraph with
# placeholder text, followed by two empty paragraphs. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq "Game Mechanics:") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $idx = $target.Index
    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = "IK HEB ECHT GEEN IDEE WAT IK HIER NEER ZOU MOETEN ZETTEN SIMON DOE JIJ DIT MAAR."

    $newPara2 = $d.Paragraphs.Item($idx + 1)
    $newPara2.Range.InsertParagraphAfter()
    $d.Paragraphs.Item($idx + 2).Range.InsertParagraphAfter()
}

